$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("begroting")
$ws2 = $wb.Worksheets.Item("Sheet1")

# --- begroting (sheet1) edits ---

# Row 26 now holds two labels: "week number" in A26 (new) and the
# previous "hours spent" label shifted to B26.
$ws1.Cells.Item(26, 1).Value = "week number"
$ws1.Cells.Item(26, 2).Value = "hours spent"

# New weeks of data appended to the table (weeks 15 and 16).
$ws1.Cells.Item(38, 1).Value = 15
$ws1.Cells.Item(38, 2).Value = 34
$ws1.Cells.Item(39, 1).Value = 16

# New "budget" label above the totals row.
$ws1.Cells.Item(41, 3).Value = "budget"

# Highlight the "Grid refinement from samples" row using the built-in
# "Neutral" cell style (used-where averaging was added).
$ws1.Range("B19").Style = "Neutral"

# --- Sheet1 (the dx / sample sheet) edits ---
$ws2.Range("B6").Value = 5

# --- View / selection state ---
[void]$ws2.Range("B9").Select()
[void]$ws1.Range("B36").Select()
[void]$ws1.Activate()
